$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1. New "age group" header in A1 (existing B1:G1 headers unchanged)
# ---------------------------------------------------------------
$ws.Range("A1").Value = "age group"

# ---------------------------------------------------------------
# 2. Column width tweaks (col A wider, col F wider)
# ---------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 20.01
$ws.Columns.Item(6).ColumnWidth = 15.97

# ---------------------------------------------------------------
# 3. New block of cells, rows 25-35 and row 38
#    (second, hospitalization-fraction table added below the
#    existing asymptomatic-fraction table)
# ---------------------------------------------------------------

# --- header row 25 ---
$ws.Range("A25").Value = "age group"
$ws.Range("B25").Value = "n0 individuals"
$ws.Range("C25").Value = "fraction"
$ws.Range("D25").Value = "fitted h"
$ws.Range("E25").Value = "multiplier"
$ws.Range("F25").Value = "pop weighted h"
$ws.Range("G25").Value = "desired h"
$ws.Range("H25").Value = "rounded h (%)"

# --- age-group labels (column A) ---
$ws.Range("A26").Value = "0-10"
$ws.Range("A27").Value = "10-20"
$ws.Range("A28").Value = "20-30"
$ws.Range("A29").Value = "30-40"
$ws.Range("A30").Value = "40-50"
$ws.Range("A31").Value = "50-60"
$ws.Range("A32").Value = "60-70"
$ws.Range("A33").Value = "70-80"
$ws.Range("A34").Value = "80-inf"
$ws.Range("A35").Value = "total"

# --- n0 individuals (column B) ---
$ws.Range("B26").Value = 1305219
$ws.Range("B27").Value = 1298970
$ws.Range("B28").Value = 1395385
$ws.Range("B29").Value = 1498535
$ws.Range("B30").Value = 1524152
$ws.Range("B31").Value = 1601891
$ws.Range("B32").Value = 1347696
$ws.Range("B33").Value = 908725
$ws.Range("B34").Value = 658753
$ws.Range("B35").Formula = '=SUM(B26:B34)'

# --- fraction (column C) ---
$ws.Range("C26").Formula = '=B26/$B$35'
$ws.Range("C27").Formula = '=B27/$B$35'
$ws.Range("C28").Formula = '=B28/$B$35'
$ws.Range("C29").Formula = '=B29/$B$35'
$ws.Range("C30").Formula = '=B30/$B$35'
$ws.Range("C31").Formula = '=B31/$B$35'
$ws.Range("C32").Formula = '=B32/$B$35'
$ws.Range("C33").Formula = '=B33/$B$35'
$ws.Range("C34").Formula = '=B34/$B$35'
$ws.Range("C35").Formula = '=B35/$B$35'

# --- fitted h (column D) ---
$ws.Range("D26").Value = 0.015
$ws.Range("D27").Value = 0.02
$ws.Range("D28").Value = 0.03
$ws.Range("D29").Value = 0.03
$ws.Range("D30").Value = 0.03
$ws.Range("D31").Value = 0.06
$ws.Range("D32").Value = 0.14
$ws.Range("D33").Value = 0.3
$ws.Range("D34").Value = 0.76

# --- multiplier (column E, only E26 populated) ---
$ws.Range("E26").Value = 0.746784953353961

# --- pop weighted h * fraction (column F) ---
$ws.Range("F26").Formula = '=D26*$E$26*C26'
$ws.Range("F27").Formula = '=D27*$E$26*C27'
$ws.Range("F28").Formula = '=D28*$E$26*C28'
$ws.Range("F29").Formula = '=D29*$E$26*C29'
$ws.Range("F30").Formula = '=D30*$E$26*C30'
$ws.Range("F31").Formula = '=D31*$E$26*C31'
$ws.Range("F32").Formula = '=D32*$E$26*C32'
$ws.Range("F33").Formula = '=D33*$E$26*C33'
$ws.Range("F34").Formula = '=D34*$E$26*C34'
$ws.Range("F35").Formula = '=SUM(F26:F34)'
$ws.Range("F35").Font.Bold = $true

# --- desired h (column G) ---
$ws.Range("G26").Formula = '=D26*$E$26'
$ws.Range("G27").Formula = '=D27*$E$26'
$ws.Range("G28").Formula = '=D28*$E$26'
$ws.Range("G29").Formula = '=D29*$E$26'
$ws.Range("G30").Formula = '=D30*$E$26'
$ws.Range("G31").Formula = '=D31*$E$26'
$ws.Range("G32").Formula = '=D32*$E$26'
$ws.Range("G33").Formula = '=D33*$E$26'
$ws.Range("G34").Formula = '=D34*$E$26'

# --- rounded h (%) (column H) ---
$ws.Range("H26").Formula = '=ROUND(G26,3)*100'
$ws.Range("H27").Formula = '=ROUND(G27,3)*100'
$ws.Range("H28").Formula = '=ROUND(G28,3)*100'
$ws.Range("H29").Formula = '=ROUND(G29,3)*100'
$ws.Range("H30").Formula = '=ROUND(G30,3)*100'
$ws.Range("H31").Formula = '=ROUND(G31,3)*100'
$ws.Range("H32").Formula = '=ROUND(G32,3)*100'
$ws.Range("H33").Formula = '=ROUND(G33,3)*100'
$ws.Range("H34").Formula = '=ROUND(G34,3)*100'

# --- desired pop weighted h, row 38 ---
$ws.Range("A38").Value = "desired pop weighted h"
$ws.Range("B38").Value = 0.08
$ws.Range("B38").Font.Bold = $true

# ---------------------------------------------------------------
# 4. Selection moves to F19
# ---------------------------------------------------------------
$ws.Range("F19").Select()
